# Update the "想去人数" (interest count) values in column F across the
# workbook's sheets, reflecting the refreshed data snapshot generated at
# commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1411
$ws1.Range("F6").Value  = 669
$ws1.Range("F12").Value = 31832
$ws1.Range("F13").Value = 6867
$ws1.Range("F15").Value = 344
$ws1.Range("F19").Value = 87
$ws1.Range("F20").Value = 42
$ws1.Range("F22").Value = 92
$ws1.Range("F25").Value = 376
$ws1.Range("F26").Value = 424
$ws1.Range("F28").Value = 184
$ws1.Range("F31").Value = 282
$ws1.Range("F33").Value = 715
$ws1.Range("F34").Value = 106
$ws1.Range("F37").Value = 282
$ws1.Range("F38").Value = 50

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 4
$ws2.Range("F5").Value = 145

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1424

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1424
$ws4.Range("F6").Value  = 1411
$ws4.Range("F8").Value  = 669
$ws4.Range("F14").Value = 145
$ws4.Range("F20").Value = 6867
$ws4.Range("F22").Value = 344
$ws4.Range("F27").Value = 87
$ws4.Range("F29").Value = 42
$ws4.Range("F32").Value = 92
$ws4.Range("F35").Value = 376
$ws4.Range("F36").Value = 424
$ws4.Range("F38").Value = 184
$ws4.Range("F42").Value = 282
$ws4.Range("F44").Value = 106
$ws4.Range("F46").Value = 282
$ws4.Range("F47").Value = 50
